$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows (row 2 and row 3) with the new environment/URL + IDs
$ws.Range("B2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "i-preproducciongestion.segurossura.com.ar"

$ws.Range("L2").Value = "'03/05/2021"
$ws.Range("L3").Value = "'03/05/2021"

$ws.Range("U2").Value = "RGM004"
$ws.Range("V2").Value = "ABCD0RRGM004"
$ws.Range("W2").Value = "ZXC0987RGM004"

$ws.Range("U3").Value = "RGM005"
$ws.Range("V3").Value = "ABCD0RRGM005"
$ws.Range("W3").Value = "ZXC0987RGM005"

# Update the view: scroll back to default (no topLeftCell override) and select A2
$ws.Range("A2").Select()

$wb.Save()
